$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "First_Name"
$ws.Range("C1").Value = "Last_Name"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Patrick"
$ws.Range("C2").Value = "Yex"

$ws.Range("B4").Select()
